$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters: A=rank, B=Coin, C=Link, D=Price, E=Volume(1h)
# Price values (column D) can look numeric (e.g. "1.00", "595.02"), so they are
# written with a leading apostrophe to force text, then the cell style is reset
# back to "Normal" so no stray number-format style gets attached to the cell.
# Volume values (column E) already contain spaces/percent signs so Excel keeps
# them as text automatically.

function Set-Price($row, $value) {
    $c = $ws.Cells.Item($row, 4)
    $c.Value = "'" + $value
    $c.Style = "Normal"
}

function Set-Volume($row, $value) {
    $ws.Cells.Item($row, 5).Value = $value
}

Set-Price  2  "66.042.22"
Set-Volume 2  "  +1.37%  "

Set-Price  3  "3.199.82"
Set-Volume 3  "  +1.19%  "

Set-Price  4  "1.00"
Set-Volume 4  "  +0.10%  "

Set-Price  5  "595.02"
Set-Volume 5  "  +3.14%  "

Set-Price  6  "153.93"
Set-Volume 6  "  +2.41%  "

Set-Price  7  "1.00"
Set-Volume 7  "  +0.00%  "

Set-Price  8  "3.196.65"
Set-Volume 8  "  +1.21%  "

Set-Volume 9  "  +1.28%  "

Set-Volume 10 "  +0.22%  "

Set-Volume 11 "  -0.79%  "

Set-Price  12 "0.514"
Set-Volume 12 "  +2.55%  "

Set-Volume 13 "  +2.17%  "

Set-Price  14 "39.01"
Set-Volume 14 "  +4.61%  "

Set-Price  15 "3.729.48"
Set-Volume 15 "  +1.45%  "

Set-Price  16 "66.028.75"
Set-Volume 16 "  +1.37%  "

Set-Price  17 "7.43"
Set-Volume 17 "  +4.25%  "

Set-Price  18 "3.210.61"
Set-Volume 18 "  +1.71%  "

Set-Volume 19 "  +0.33%  "

Set-Price  20 "510.50"
Set-Volume 20 "  +0.09%  "

Set-Volume 21 "  +3.44%  "

Set-Price  22 "0.741"
Set-Volume 22 "  +3.02%  "

Set-Price  23 "15.23"
Set-Volume 23 "  -0.70%  "

Set-Price  24 "7.99"
Set-Volume 24 "  +2.87%  "

Set-Price  25 "84.94"
Set-Volume 25 "  +0.46%  "

Set-Volume 26 "  +0.01%  "

Set-Price  27 "9.40"
Set-Volume 27 "  +5.39%  "

Set-Price  28 "2.99"
Set-Volume 28 "  +2.74%  "

Set-Price  29 "2.27"
Set-Volume 29 "  +4.11%  "

Set-Price  30 "6.89"
Set-Volume 30 "  +9.62%  "

Set-Price  31 "2.87"
Set-Volume 31 "  +3.33%  "

Set-Price  32 "28.28"
Set-Volume 32 "  +1.79%  "

Set-Volume 33 "  +3.04%  "

Set-Price  34 "1.00"
Set-Volume 34 "  +0.32%  "

Set-Price  35 "6.56"
Set-Volume 35 "  +0.09%  "

Set-Price  36 "55.03"
Set-Volume 36 "  -0.68%  "

Set-Price  37 "0.0904"
Set-Volume 37 "  +0.16%  "

Set-Price  38 "485.85"
Set-Volume 38 "  +3.71%  "

Set-Price  39 "0.0419"
Set-Volume 39 "  -1.96%  "

Set-Price  40 "2.93"
Set-Volume 40 "  -3.74%  "

Set-Price  41 "8.83"
Set-Volume 41 "  +1.77%  "

# Rows 42 and 43 swap content entirely (TheGraph <-> Kaspa)
$ws.Cells.Item(42, 2).Value = "Kaspa"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-Price  42 "0.122"
Set-Volume 42 "  +3.58%  "

$ws.Cells.Item(43, 2).Value = "TheGraph"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-Price  43 "0.299"
Set-Volume 43 "  +5.46%  "

Set-Price  44 "0.0₃0649"
Set-Volume 44 "  +9.58%  "

Set-Price  45 "2.938.74"
Set-Volume 45 "  -3.89%  "

Set-Price  46 "2.42"
Set-Volume 46 "  -1.34%  "

Set-Price  47 "28.42"
Set-Volume 47 "  -1.43%  "

Set-Volume 48 "  +0.00%  "

Set-Price  49 "0.116"
Set-Volume 49 "  +1.48%  "

Set-Volume 50 "  +1.75%  "

Set-Price  51 "2.58"
Set-Volume 51 "  +4.31%  "
